$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 496.2
$ws.Range("J6").Value = 1300
$ws.Range("L6").Value = 3900
$ws.Range("N6").Value = -4124
$ws.Range("H40").Value = 1613.7142
$ws.Range("J40").Value = 1498
$ws.Range("L40").Value = 1498
$ws.Range("N40").Value = -1848
$ws.Range("H45").Value = 13004.25
$ws.Range("J45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("N45").Value = $null
$ws.Range("H51").Value = 3485.3572
$ws.Range("I51").Value = 2897.8
$ws.Range("J51").Value = 3811.7778
$ws.Range("K51").Value = 2897.8
$ws.Range("L51").Value = 3811.7778
$ws.Range("M51").Value = -2413.8
$ws.Range("N51").Value = -4779.7778
$ws.Range("H62").Value = 1913.619
$ws.Range("I62").Value = 2040.3529
$ws.Range("J62").Value = 1375
$ws.Range("K62").Value = 2040.3529
$ws.Range("L62").Value = 1375
$ws.Range("M62").Value = -1416.3529
$ws.Range("N62").Value = -2623
$ws.Range("H65").Value = 1913.619
$ws.Range("I65").Value = 2040.3529
$ws.Range("J65").Value = 1375
$ws.Range("K65").Value = 10201.7645
$ws.Range("L65").Value = 6875
$ws.Range("M65").Value = -7081.764500000001
$ws.Range("N65").Value = -13115
$ws.Range("H74").Value = 5278.3
$ws.Range("I74").Value = 5195.75
$ws.Range("J74").Value = 5333.3335
$ws.Range("K74").Value = 5195.75
$ws.Range("L74").Value = 5333.3335
$ws.Range("M74").Value = -4259.75
$ws.Range("N74").Value = -7205.3335
$ws.Range("H77").Value = 5278.3
$ws.Range("I77").Value = 5195.75
$ws.Range("J77").Value = 5333.3335
$ws.Range("K77").Value = 25978.75
$ws.Range("L77").Value = 26666.6675
$ws.Range("M77").Value = -21298.75
$ws.Range("N77").Value = -36026.6675
$ws.Range("H88").Value = 2101.4167
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 2101.4167
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 2101.4167
$ws.Range("M88").Value = $null
$ws.Range("N88").Value = -2913.4167
$ws.Range("H91").Value = 2101.4167
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 2101.4167
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 2101.4167
$ws.Range("M91").Value = $null
$ws.Range("N91").Value = -4909.4167
$ws.Range("H98").Value = 4700.9614
$ws.Range("I98").Value = 2884.375
$ws.Range("K98").Value = 2884.375
$ws.Range("M98").Value = -1386.375
$ws.Range("H112").Value = 1931.3823
$ws.Range("I112").Value = 822
$ws.Range("J112").Value = 2079.3
$ws.Range("K112").Value = 2466
$ws.Range("L112").Value = 6237.900000000001
$ws.Range("M112").Value = -1358
$ws.Range("N112").Value = -8453.900000000001
$ws.Range("H121").Value = 1937.2222
$ws.Range("J121").Value = 2234.2856
$ws.Range("L121").Value = 6702.8568
$ws.Range("N121").Value = -10196.8568
$ws.Range("H122").Value = 4700.9614
$ws.Range("I122").Value = 2884.375
$ws.Range("K122").Value = 8653.125
$ws.Range("M122").Value = -6203.125
$ws.Range("H125").Value = 504899.5
$ws.Range("J125").Value = 9800
$ws.Range("L125").Value = 88200
$ws.Range("N125").Value = -93120
$ws.Range("H129").Value = 1204.3208
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 1204.3208
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 3612.9624
$ws.Range("M129").Value = $null
$ws.Range("N129").Value = -13612.9624

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H62").Value = 30000
$ws.Range("I62").Value = 30000
$ws.Range("K62").Value = 30000
$ws.Range("M62").Value = -29376
$ws.Range("H65").Value = 30000
$ws.Range("I65").Value = 30000
$ws.Range("K65").Value = 90000
$ws.Range("M65").Value = -86880
$ws.Range("H123").Value = 24129
$ws.Range("J123").Value = 24129
$ws.Range("L123").Value = 24129
$ws.Range("N123").Value = -33929

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 145566.14
$ws.Range("I86").Value = 3159.3333
$ws.Range("K86").Value = 3159.3333
$ws.Range("M86").Value = -2036.3333
$ws.Range("H89").Value = 145566.14
$ws.Range("I89").Value = 3159.3333
$ws.Range("K89").Value = 15796.6665
$ws.Range("M89").Value = -10180.6665
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").Value = $null
$ws.Range("H134").Value = 3152.6296
$ws.Range("I134").Value = 3114.1333
$ws.Range("J134").Value = 3200.75
$ws.Range("K134").Value = 9342.3999
$ws.Range("L134").Value = 9602.25
$ws.Range("M134").Value = -6807.3999
$ws.Range("N134").Value = -14672.25
$ws.Range("H135").Value = 74640
$ws.Range("J135").Value = 74640
$ws.Range("L135").Value = 74640
$ws.Range("N135").Value = -84780

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1454.5
$ws.Range("I16").Value = 1418.125
$ws.Range("J16").Value = 1600
$ws.Range("K16").Value = 1418.125
$ws.Range("L16").Value = 1600
$ws.Range("M16").Value = -1131.125
$ws.Range("N16").Value = -2174
$ws.Range("H22").Value = 362.96155
$ws.Range("I22").Value = 207.7619
$ws.Range("J22").Value = 1014.8
$ws.Range("K22").Value = 207.7619
$ws.Range("L22").Value = 1014.8
$ws.Range("M22").Value = 142.2381
$ws.Range("N22").Value = -1714.8
$ws.Range("H60").Value = 10902.454
$ws.Range("J60").Value = 11103
$ws.Range("L60").Value = 11103
$ws.Range("N60").Value = -12125
$ws.Range("H87").Value = 39997.5
$ws.Range("J87").Value = 39997.5
$ws.Range("L87").Value = 39997.5
$ws.Range("N87").Value = -42369.5
$ws.Range("H90").Value = 39997.5
$ws.Range("J90").Value = 39997.5
$ws.Range("L90").Value = 119992.5
$ws.Range("N90").Value = -131848.5
$ws.Range("H105").Value = 2262
$ws.Range("I105").Value = 2070
$ws.Range("J105").Value = 2550
$ws.Range("K105").Value = 2070
$ws.Range("L105").Value = 2550
$ws.Range("M105").Value = -323
$ws.Range("N105").Value = -6044
$ws.Range("H113").Value = 1454.5
$ws.Range("I113").Value = 1418.125
$ws.Range("J113").Value = 1600
$ws.Range("K113").Value = 1418.125
$ws.Range("L113").Value = 1600
$ws.Range("M113").Value = 751.875
$ws.Range("N113").Value = -5940
$ws.Range("H132").Value = 437846.88
$ws.Range("I132").Value = 564519.4399999999
$ws.Range("J132").Value = 3540.8572
$ws.Range("K132").Value = 1693558.32
$ws.Range("L132").Value = 10622.5716
$ws.Range("M132").Value = -1691028.32
$ws.Range("N132").Value = -15682.5716

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 799.2222
$ws.Range("I107").Value = 900.1667
$ws.Range("J107").Value = 597.3333
$ws.Range("K107").Value = 2700.5001
$ws.Range("L107").Value = 1791.9999
$ws.Range("M107").Value = -780.5001000000002
$ws.Range("N107").Value = -5631.9999
$ws.Range("H113").Value = 657.3570999999999
$ws.Range("J113").Value = 778.44446
$ws.Range("L113").Value = 2335.33338
$ws.Range("N113").Value = -6675.33338
$ws.Range("H132").Value = 1697.6111
$ws.Range("I132").Value = 1309
$ws.Range("K132").Value = 11781
$ws.Range("M132").Value = -9251

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2905.7144
$ws.Range("I80").Value = 2642.8572
$ws.Range("J80").Value = 3168.5715
$ws.Range("K80").Value = 2642.8572
$ws.Range("L80").Value = 3168.5715
$ws.Range("M80").Value = -1644.8572
$ws.Range("N80").Value = -5164.5715
$ws.Range("H83").Value = 2905.7144
$ws.Range("I83").Value = 2642.8572
$ws.Range("J83").Value = 3168.5715
$ws.Range("K83").Value = 13214.286
$ws.Range("L83").Value = 15842.8575
$ws.Range("M83").Value = -8222.286
$ws.Range("N83").Value = -25826.8575
$ws.Range("H93").Value = 31041.666
$ws.Range("J93").Value = 31041.666
$ws.Range("L93").Value = 31041.666
$ws.Range("N93").Value = -34785.666
$ws.Range("H109").Value = 9223.0625
$ws.Range("J109").Value = 9223.0625
$ws.Range("L109").Value = 9223.0625
$ws.Range("N109").Value = -11303.0625
$ws.Range("H132").Value = 2480.4666
$ws.Range("I132").Value = 1863.875
$ws.Range("J132").Value = 3185.1428
$ws.Range("K132").Value = 5591.625
$ws.Range("L132").Value = 9555.428400000001
$ws.Range("M132").Value = -3061.625
$ws.Range("N132").Value = -14615.4284
$ws.Range("H135").Value = 69995
$ws.Range("J135").Value = 69995
$ws.Range("L135").Value = 69995
$ws.Range("N135").Value = -80135

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H139").Value = 39933.332
$ws.Range("J139").Value = 39933.332
$ws.Range("L139").Value = 39933.332
$ws.Range("N139").Value = -50213.332

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 35000
$ws.Range("J64").Value = 35000
$ws.Range("L64").Value = 35000
$ws.Range("N64").Value = -35496
$ws.Range("H67").Value = 35000
$ws.Range("J67").Value = 35000
$ws.Range("L67").Value = 35000
$ws.Range("N67").Value = -36716

